$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value = 44195
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 20
$ws.Cells.Item(3, 14).Value = 15000
$ws.Cells.Item(3, 15).Value = 15000
$ws.Cells.Item(3, 16).Value = 15000
$ws.Cells.Item(3, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(3, 19).Value = 1500
$ws.Cells.Item(3, 20).Value = 10

# Row 4
$ws.Cells.Item(4, 4).Value = 44391
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 15
$ws.Cells.Item(4, 14).Value = 1500
$ws.Cells.Item(4, 15).Value = 1500
$ws.Cells.Item(4, 16).Value = 1500
$ws.Cells.Item(4, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(4, 19).Value = 1500
$ws.Cells.Item(4, 20).Value = 1

# Row 5
$ws.Cells.Item(5, 4).Value = 44391
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 20
$ws.Cells.Item(5, 14).Value = 1000
$ws.Cells.Item(5, 15).Value = 1000
$ws.Cells.Item(5, 16).Value = 1000
$ws.Cells.Item(5, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(5, 19).Value = 1000
$ws.Cells.Item(5, 20).Value = 1

# Row 6
$ws.Cells.Item(6, 4).Value = 44904
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 45
$ws.Cells.Item(6, 14).Value = 15000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 15000
$ws.Cells.Item(6, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(6, 19).Value = 1500
$ws.Cells.Item(6, 20).Value = 10

# Row 7
$ws.Cells.Item(7, 4).Value = 44904
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 60
$ws.Cells.Item(7, 14).Value = 10000
$ws.Cells.Item(7, 15).Value = 10000
$ws.Cells.Item(7, 16).Value = 10000
$ws.Cells.Item(7, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(7, 19).Value = 1000
$ws.Cells.Item(7, 20).Value = 10

# Row 8
$ws.Cells.Item(8, 4).Value = 44371
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 20
$ws.Cells.Item(8, 14).Value = 1800
$ws.Cells.Item(8, 15).Value = 1800
$ws.Cells.Item(8, 16).Value = 1800
$ws.Cells.Item(8, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(8, 19).Value = 1800
$ws.Cells.Item(8, 20).Value = 1

# Row 9
$ws.Cells.Item(9, 4).Value = 44371
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 1200
$ws.Cells.Item(9, 15).Value = 1200
$ws.Cells.Item(9, 16).Value = 1200
$ws.Cells.Item(9, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(9, 19).Value = 1200
$ws.Cells.Item(9, 20).Value = 1

# Row 10
$ws.Cells.Item(10, 4).Value = 44336
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 10
$ws.Cells.Item(10, 14).Value = 1500
$ws.Cells.Item(10, 15).Value = 1500
$ws.Cells.Item(10, 16).Value = 1500
$ws.Cells.Item(10, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(10, 19).Value = 1500
$ws.Cells.Item(10, 20).Value = 1

# Row 11
$ws.Cells.Item(11, 4).Value = 44343
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 20
$ws.Cells.Item(11, 14).Value = 1700
$ws.Cells.Item(11, 15).Value = 1700
$ws.Cells.Item(11, 16).Value = 1700
$ws.Cells.Item(11, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(11, 19).Value = 1700
$ws.Cells.Item(11, 20).Value = 1

# Row 12
$ws.Cells.Item(12, 4).Value = 44400
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 25
$ws.Cells.Item(12, 14).Value = 1500
$ws.Cells.Item(12, 15).Value = 1500
$ws.Cells.Item(12, 16).Value = 1500
$ws.Cells.Item(12, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(12, 19).Value = 1500
$ws.Cells.Item(12, 20).Value = 1

# Row 13
$ws.Cells.Item(13, 4).Value = 44880
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 20000
$ws.Cells.Item(13, 15).Value = 20000
$ws.Cells.Item(13, 16).Value = 20000
$ws.Cells.Item(13, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(13, 19).Value = 2000
$ws.Cells.Item(13, 20).Value = 10

# Row 14
$ws.Cells.Item(14, 4).Value = 44880
$ws.Cells.Item(14, 12).Value = 'Segunda'
$ws.Cells.Item(14, 13).Value = 180
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 15000
$ws.Cells.Item(14, 16).Value = 15000
$ws.Cells.Item(14, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(14, 19).Value = 1500
$ws.Cells.Item(14, 20).Value = 10

# Row 15
$ws.Cells.Item(15, 4).Value = 44309
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 10
$ws.Cells.Item(15, 14).Value = 1600
$ws.Cells.Item(15, 15).Value = 1600
$ws.Cells.Item(15, 16).Value = 1600
$ws.Cells.Item(15, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(15, 19).Value = 1600
$ws.Cells.Item(15, 20).Value = 1
